# Update "want to go" counts (column F) in sheets "展览" (Exhibition) and "全部类型" (All types)
$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Cells.Item(3, 6).Value = 274
$wsExhibition.Cells.Item(4, 6).Value = 1811
$wsExhibition.Cells.Item(6, 6).Value = 557
$wsExhibition.Cells.Item(7, 6).Value = 5172
$wsExhibition.Cells.Item(11, 6).Value = 999
$wsExhibition.Cells.Item(12, 6).Value = 355
$wsExhibition.Cells.Item(13, 6).Value = 1318
$wsExhibition.Cells.Item(16, 6).Value = 3039
$wsExhibition.Cells.Item(17, 6).Value = 1874
$wsExhibition.Cells.Item(21, 6).Value = 112
$wsExhibition.Cells.Item(22, 6).Value = 652
$wsExhibition.Cells.Item(24, 6).Value = 333
$wsExhibition.Cells.Item(25, 6).Value = 43
$wsExhibition.Cells.Item(26, 6).Value = 3424
$wsExhibition.Cells.Item(27, 6).Value = 1066
$wsExhibition.Cells.Item(28, 6).Value = 2698
$wsExhibition.Cells.Item(29, 6).Value = 272
$wsExhibition.Cells.Item(30, 6).Value = 1691
$wsExhibition.Cells.Item(31, 6).Value = 3883
$wsExhibition.Cells.Item(35, 6).Value = 1230
$wsExhibition.Cells.Item(37, 6).Value = 968
$wsExhibition.Cells.Item(38, 6).Value = 1229
$wsExhibition.Cells.Item(39, 6).Value = 45
$wsExhibition.Cells.Item(40, 6).Value = 975
$wsExhibition.Cells.Item(41, 6).Value = 633
$wsExhibition.Cells.Item(42, 6).Value = 455
$wsExhibition.Cells.Item(43, 6).Value = 391
$wsExhibition.Cells.Item(44, 6).Value = 303
$wsExhibition.Cells.Item(45, 6).Value = 3531

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Cells.Item(3, 6).Value = 274
$wsAllTypes.Cells.Item(4, 6).Value = 1811
$wsAllTypes.Cells.Item(6, 6).Value = 557
$wsAllTypes.Cells.Item(7, 6).Value = 5172
$wsAllTypes.Cells.Item(12, 6).Value = 355
$wsAllTypes.Cells.Item(13, 6).Value = 1318
$wsAllTypes.Cells.Item(14, 6).Value = 3039
$wsAllTypes.Cells.Item(16, 6).Value = 1874
$wsAllTypes.Cells.Item(23, 6).Value = 112
$wsAllTypes.Cells.Item(26, 6).Value = 333
$wsAllTypes.Cells.Item(27, 6).Value = 3424
$wsAllTypes.Cells.Item(30, 6).Value = 1066
$wsAllTypes.Cells.Item(31, 6).Value = 2698
$wsAllTypes.Cells.Item(32, 6).Value = 1691
$wsAllTypes.Cells.Item(33, 6).Value = 3883
$wsAllTypes.Cells.Item(37, 6).Value = 1230
$wsAllTypes.Cells.Item(39, 6).Value = 968
$wsAllTypes.Cells.Item(41, 6).Value = 1229
$wsAllTypes.Cells.Item(42, 6).Value = 45
$wsAllTypes.Cells.Item(43, 6).Value = 975
$wsAllTypes.Cells.Item(44, 6).Value = 633
$wsAllTypes.Cells.Item(45, 6).Value = 391
$wsAllTypes.Cells.Item(48, 6).Value = 303
$wsAllTypes.Cells.Item(49, 6).Value = 3531

$wb.Save()
